$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaText = "Meta description: Read our review of Chang Thai slot game and play for free. Find out what it offers, including bonus game, multiple paylines, and Autoplay feature."
$metaPara.Range.Text = $metaText

# Bold just the "Meta description" label at the start of the new paragraph.
$labelStart = $metaPara.Range.Start
$labelEnd = $labelStart + ("Meta description").Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Bold = 1

# ------------------------------------------------------------------
# 2) Near the end of the document there used to be two paragraphs:
#      - a bold duplicate of the page title
#      - an italic "meta description" style blurb
#    Remove the bold duplicate-title paragraph entirely, and replace
#    the text of the remaining italic paragraph with the new image
#    generation prompt (keeping its italic run formatting intact).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

$count2 = $d.Paragraphs.Count
$imgPromptPara = $d.Paragraphs($count2)
$imgStart = $imgPromptPara.Range.Start
$imgEnd = $imgPromptPara.Range.End
$imgRange = $d.Range($imgStart, $imgEnd)

$imgText = 'Please create a cartoon-style image featuring a happy Maya warrior with glasses for the game "Chang Thai". The image should be vibrant and eye-catching, with the Maya warrior as the main focus. The warrior should be shown wearing glasses that reflect the lush vegetation of the Thai forest, and perhaps holding a golden flower with a red gem in the center to represent the game''s Scatter symbol. The background of the image should feature the Thai temple and animated elephant, as mentioned in the game review, to tie it back to the game''s setting. The overall tone of the image should be fun and adventurous, inviting players to join the journey through the southern Thai forest.'

$imgRange.Text = $imgText
